$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.878.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.394.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "504.64"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.89"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.552"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.399.26"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0976"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.63"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.819.53"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.812.64"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.72"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.443.24"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.23"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.07"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "310.25"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.60"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.52"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.151"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.41"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.66"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0723"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.65"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.02%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.85"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.59%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.93"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.82"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.827"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.86"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.84%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "131.41"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.37"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.84"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.568"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0914"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "250.87"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0487"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.05"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.21%  "
